$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows after row 410 (new weekly data), shifting old rows 411-431 down to 415-435
$ws.Rows("411:414").Insert()

# Row 411
$ws.Cells.Item(411,1).Value = 10
$ws.Cells.Item(411,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(411,3).Value = 'La Araucanía'
$ws.Cells.Item(411,4).Value = 44714
$ws.Cells.Item(411,5).Value = 9
$ws.Cells.Item(411,6).Value = 'Fruta'
$ws.Cells.Item(411,7).Value = 100101
$ws.Cells.Item(411,8).Value = 'Berries'
$ws.Cells.Item(411,9).Value = 100101007
$ws.Cells.Item(411,10).Value = 'Kiwi'
$ws.Cells.Item(411,11).Value = 'Hayward'
$ws.Cells.Item(411,12).Value = 'Primera'
$ws.Cells.Item(411,13).Value = 150
$ws.Cells.Item(411,14).Value = 13000
$ws.Cells.Item(411,15).Value = 13000
$ws.Cells.Item(411,16).Value = 13000
$ws.Cells.Item(411,17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(411,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(411,19).Value = 722
$ws.Cells.Item(411,20).Value = 18

# Row 412
$ws.Cells.Item(412,1).Value = 10
$ws.Cells.Item(412,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(412,3).Value = 'La Araucanía'
$ws.Cells.Item(412,4).Value = 44714
$ws.Cells.Item(412,5).Value = 9
$ws.Cells.Item(412,6).Value = 'Fruta'
$ws.Cells.Item(412,7).Value = 100101
$ws.Cells.Item(412,8).Value = 'Berries'
$ws.Cells.Item(412,9).Value = 100101007
$ws.Cells.Item(412,10).Value = 'Kiwi'
$ws.Cells.Item(412,11).Value = 'Hayward'
$ws.Cells.Item(412,12).Value = 'Primera'
$ws.Cells.Item(412,13).Value = 6
$ws.Cells.Item(412,14).Value = 350000
$ws.Cells.Item(412,15).Value = 350000
$ws.Cells.Item(412,16).Value = 350000
$ws.Cells.Item(412,17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(412,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(412,19).Value = 778
$ws.Cells.Item(412,20).Value = 450

# Row 413
$ws.Cells.Item(413,1).Value = 10
$ws.Cells.Item(413,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(413,3).Value = 'La Araucanía'
$ws.Cells.Item(413,4).Value = 44714
$ws.Cells.Item(413,5).Value = 9
$ws.Cells.Item(413,6).Value = 'Fruta'
$ws.Cells.Item(413,7).Value = 100101
$ws.Cells.Item(413,8).Value = 'Berries'
$ws.Cells.Item(413,9).Value = 100101007
$ws.Cells.Item(413,10).Value = 'Kiwi'
$ws.Cells.Item(413,11).Value = 'Hayward'
$ws.Cells.Item(413,12).Value = 'Tercera'
$ws.Cells.Item(413,13).Value = 40
$ws.Cells.Item(413,14).Value = 8000
$ws.Cells.Item(413,15).Value = 8000
$ws.Cells.Item(413,16).Value = 8000
$ws.Cells.Item(413,17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(413,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(413,19).Value = 444
$ws.Cells.Item(413,20).Value = 18

# Row 414
$ws.Cells.Item(414,1).Value = 10
$ws.Cells.Item(414,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(414,3).Value = 'La Araucanía'
$ws.Cells.Item(414,4).Value = 44714
$ws.Cells.Item(414,5).Value = 9
$ws.Cells.Item(414,6).Value = 'Fruta'
$ws.Cells.Item(414,7).Value = 100101
$ws.Cells.Item(414,8).Value = 'Berries'
$ws.Cells.Item(414,9).Value = 100101007
$ws.Cells.Item(414,10).Value = 'Kiwi'
$ws.Cells.Item(414,11).Value = 'Hayward'
$ws.Cells.Item(414,12).Value = 'Tercera'
$ws.Cells.Item(414,13).Value = 3
$ws.Cells.Item(414,14).Value = 150000
$ws.Cells.Item(414,15).Value = 150000
$ws.Cells.Item(414,16).Value = 150000
$ws.Cells.Item(414,17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(414,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(414,19).Value = 333
$ws.Cells.Item(414,20).Value = 450
